$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.47000000000054
$ws.Range("G2").Value = [double]"8.412080368280783e-06"
$ws.Range("H2").Value = [double]"3.817782814854226e-05"
$ws.Range("K2").Value = 5.765171682625756
$ws.Range("L2").Value = "[2.628374463533204, 8.901968901718307]"
$ws.Range("M2").Value = 0.0003462748986908792
$ws.Range("N2").Value = 0.0003462748986908792
$ws.Range("O2").Value = -0.9056843686024632
$ws.Range("P2").Value = "[-1.4088423511593877, -0.4025263860455386]"
$ws.Range("Q2").Value = 0.0004557246923508895
$ws.Range("R2").Value = 0.0004557246923508895
$ws.Range("S2").Value = 11.44577489943677
$ws.Range("T2").Value = "[9.81141426750001, 13.080135531373537]"
$ws.Range("W2").Value = 3.671351351351433
$ws.Range("X2").Value = 1.631711711711747
$ws.Range("Y2").Value = 5.710990990991119

# Row 3 updates
$ws.Range("E3").Value = 23.99000000000031
$ws.Range("G3").Value = [double]"1.360274543005247e-07"
$ws.Range("H3").Value = [double]"3.503536500566756e-06"
$ws.Range("K3").Value = 5.415194213472541
$ws.Range("L3").Value = "[3.0190192907260514, 7.81136913621903]"
$ws.Range("M3").Value = [double]"1.135957846853941e-05"
$ws.Range("N3").Value = [double]"2.271915693707882e-05"
$ws.Range("O3").Value = 2.735921530153274
$ws.Range("P3").Value = "[2.257921446724196, 3.213921613582351]"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 11.20575045319496
$ws.Range("T3").Value = "[9.892643995937174, 12.518856910452737]"
$ws.Range("W3").Value = 13.54390390390408
$ws.Range("X3").Value = 11.71883883883899
$ws.Range("Y3").Value = 15.36896896896917
